$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7333541512489319
$ws.Range("B1").Value = 1.397609710693359
$ws.Range("C1").Value = 4.485628128051758
$ws.Range("D1").Value = 1.821876287460327
$ws.Range("E1").Value = 1.133834481239319
